$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts old E->F, F->G, G->H)
$ws.Columns.Item(5).Insert()

# New header for the inserted column E
$ws.Range("E1").Value = "Cost Of Investment *"

# New "Cost Of Investment" values for rows 2-7
$ws.Range("E2").Value = 100000
$ws.Range("E3").Value = 200000
$ws.Range("E4").Value = 90000
$ws.Range("E5").Value = 100000
$ws.Range("E6").Value = 200000
$ws.Range("E7").Value = 200000

# Match the number formatting used in column D (numFmtId 4, thousands + 2 decimals)
$ws.Range("E2:E7").NumberFormat = "#,##0.00"

# Give the new column an explicit (non-autofit) width
$ws.Columns.Item(5).ColumnWidth = 18.6875

# Update selection to match the new layout
$ws.Range("E8").Select()
